$d = $word.ActiveDocument

# Build the long separator line (60 U+2500 BOX DRAWINGS LIGHT HORIZONTAL chars)
$sepChar = [char]0x2500
$sepLine = ""
for ($i = 0; $i -lt 60; $i++) { $sepLine = $sepLine + $sepChar }

# Walk paragraphs back-to-front so deleting one doesn't shift the indices
# of the ones we still need to examine.
$n = $d.Paragraphs.Count
for ($i = $n; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $rng = $p.Range
    $t = $rng.Text

    $isDrawing = ($rng.InlineShapes.Count -gt 0)
    $isSeparator = $t.StartsWith($sepLine)
    $isSpacerEmpty = ($t.Length -eq 1) -and ($p.Format.SpaceBefore -eq 2) -and (-not $isDrawing)

    if ($isDrawing -or $isSeparator -or $isSpacerEmpty) {
        $rng.Delete()
    }
}

Write-Output "Done. Paragraphs now: $($d.Paragraphs.Count)"
